$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Save" flag, derived from the existing G (sum) column.
# Header cell H1 - mirror the style used by the other header cells (copy format from G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-68: Save = 1 when sum (G) is a "save" outing (>= 8), else 0.
$saveValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,1,1,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,1,0,0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
